$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7117
$ws.Range("K3").Value = 7367
$ws.Range("J4").Value = 1843
$ws.Range("K4").Value = 1535
$ws.Range("K6").Value = 8132
$ws.Range("J7").Value = 29310
$ws.Range("K7").Value = 24673
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K6").Value = 180
$ws.Range("K8").Value = 1606
$ws.Range("K9").Value = 114
$ws.Range("K10").Value = 142
$ws.Range("K14").Value = 120
$ws.Range("K19").Value = 720
$ws.Range("K20").Value = 605
$ws.Range("K23").Value = 251
$ws.Range("K27").Value = 232
$ws.Range("K29").Value = 1360
$ws.Range("K31").Value = 285
$ws.Range("K33").Value = 1054
$ws.Range("K34").Value = 139
$ws.Range("K35").Value = 38
$ws.Range("K42").Value = 910
$ws.Range("K43").Value = 203
$ws.Range("K47").Value = 166
$ws.Range("K48").Value = 315
$ws.Range("K52").Value = 641
$ws.Range("K54").Value = 480
$ws.Range("K57").Value = 96
$ws.Range("K63").Value = 67
$ws.Range("K65").Value = 581
$ws.Range("K67").Value = 967
$ws.Range("J68").Value = 65
$ws.Range("K69").Value = 54
$ws.Range("K74").Value = 26
$ws.Range("K76").Value = 336
$ws.Range("K78").Value = 298
$ws.Range("K79").Value = 608
$ws.Range("K83").Value = 524
$ws.Range("K85").Value = 1128
$ws.Range("K89").Value = 369
$ws.Range("K90").Value = 236
$ws.Range("K91").Value = 292
$ws.Range("K92").Value = 92
$ws.Range("K97").Value = 197
$ws.Range("J101").Value = 29310
$ws.Range("K101").Value = 24673
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 120
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 115
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 369
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 370
$ws.Range("K3").Value = 392
$ws.Range("K7").Value = 1128
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 237
$ws.Range("K7").Value = 641
$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 54
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 447
$ws.Range("K3").Value = 486
$ws.Range("K6").Value = 537
$ws.Range("K7").Value = 1606
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K6").Value = 123
$ws.Range("K7").Value = 524
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 264
$ws.Range("K6").Value = 332
$ws.Range("K7").Value = 1054
$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 217
$ws.Range("K7").Value = 581
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 86
$ws.Range("K6").Value = 113
$ws.Range("K7").Value = 285
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 269
$ws.Range("K6").Value = 271
$ws.Range("K7").Value = 967
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 264
$ws.Range("K7").Value = 480
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 383
$ws.Range("K3").Value = 485
$ws.Range("K6").Value = 396
$ws.Range("K7").Value = 1360
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K4").Value = 42
$ws.Range("K7").Value = 315
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 211
$ws.Range("K3").Value = 217
$ws.Range("K7").Value = 720
$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 75
$ws.Range("K3").Value = 64
$ws.Range("K6").Value = 169
$ws.Range("K7").Value = 336
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 70
$ws.Range("K7").Value = 180
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K6").Value = 339
$ws.Range("K7").Value = 910
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 142
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 88
$ws.Range("K3").Value = 77
$ws.Range("K7").Value = 298
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 86
$ws.Range("K7").Value = 251
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 76
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 292
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K6").Value = 157
$ws.Range("K7").Value = 608
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 208
$ws.Range("K3").Value = 195
$ws.Range("K6").Value = 166
$ws.Range("K7").Value = 605
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 139
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K2").Value = 48
$ws.Range("K7").Value = 166
$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 38
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K2").Value = 37
$ws.Range("K7").Value = 114
$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 197
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 92
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 232
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K3").Value = 69
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 236
$ws = $wb.Worksheets.Item('North Park')
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 65
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 96
$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 203
$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 26
